$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 31/12/2025 10:02:17'
$ws1.Range("A3").Value = 'Total filas: 831'
$ws1.Range("B812").Value = '10:02:06'
$ws1.Range("C812").Value = '10:04'
$ws1.Range("D812").Value = '23_HERNANDEZ'
$ws1.Range("E812").Value = 2
$ws1.Range("F812").Value = 'LP1912'
$ws1.Range("G812").Value = '31/12/2025'
$ws1.Range("B813").Value = '10:02:06'
$ws1.Range("C813").Value = '10:14'
$ws1.Range("D813").Value = '10_OLMOS'
$ws1.Range("E813").Value = 12
$ws1.Range("F813").Value = 'LP1912'
$ws1.Range("G813").Value = '31/12/2025'
$ws1.Range("B814").Value = '10:02:06'
$ws1.Range("C814").Value = '10:25'
$ws1.Range("D814").Value = '11_ETCHEVERRY'
$ws1.Range("E814").Value = 23
$ws1.Range("F814").Value = 'LP1912'
$ws1.Range("G814").Value = '31/12/2025'
$ws1.Range("B815").Value = '10:02:06'
$ws1.Range("C815").Value = '10:25'
$ws1.Range("D815").Value = '16_SANTA ANA'
$ws1.Range("E815").Value = 23
$ws1.Range("F815").Value = 'LP1912'
$ws1.Range("G815").Value = '31/12/2025'
$ws1.Range("B816").Value = '10:02:06'
$ws1.Range("C816").Value = '10:26'
$ws1.Range("D816").Value = '15X38_ABASTO'
$ws1.Range("E816").Value = 24
$ws1.Range("F816").Value = 'LP1912'
$ws1.Range("G816").Value = '31/12/2025'
$ws1.Range("B817").Value = '10:02:06'
$ws1.Range("C817").Value = '10:34'
$ws1.Range("D817").Value = '10_OLMOS'
$ws1.Range("E817").Value = 32
$ws1.Range("F817").Value = 'LP1912'
$ws1.Range("G817").Value = '31/12/2025'
$ws1.Range("B818").Value = '10:02:06'
$ws1.Range("C818").Value = '10:34'
$ws1.Range("D818").Value = '23_HERNANDEZ'
$ws1.Range("E818").Value = 32
$ws1.Range("F818").Value = 'LP1912'
$ws1.Range("G818").Value = '31/12/2025'
$ws1.Range("B819").Value = '10:02:06'
$ws1.Range("C819").Value = '10:37'
$ws1.Range("D819").Value = '16_P MOR-SANTA ANA'
$ws1.Range("E819").Value = 35
$ws1.Range("F819").Value = 'LP1912'
$ws1.Range("G819").Value = '31/12/2025'
$ws1.Range("B820").Value = '10:02:06'
$ws1.Range("C820").Value = '10:39'
$ws1.Range("D820").Value = '15_ABASTO'
$ws1.Range("E820").Value = 37
$ws1.Range("F820").Value = 'LP1912'
$ws1.Range("G820").Value = '31/12/2025'
$ws1.Range("B821").Value = '10:02:06'
$ws1.Range("C821").Value = '10:49'
$ws1.Range("D821").Value = '16_SANTA ANA'
$ws1.Range("E821").Value = 47
$ws1.Range("F821").Value = 'LP1912'
$ws1.Range("G821").Value = '31/12/2025'
$ws1.Range("B822").Value = '10:02:06'
$ws1.Range("C822").Value = '10:51'
$ws1.Range("D822").Value = '15_ABASTO'
$ws1.Range("E822").Value = 49
$ws1.Range("F822").Value = 'LP1912'
$ws1.Range("G822").Value = '31/12/2025'
$ws1.Range("B823").Value = '10:02:06'
$ws1.Range("C823").Value = '10:54'
$ws1.Range("D823").Value = '10_OLMOS'
$ws1.Range("E823").Value = 52
$ws1.Range("F823").Value = 'LP1912'
$ws1.Range("G823").Value = '31/12/2025'
$ws1.Range("B824").Value = '10:02:06'
$ws1.Range("C824").Value = '10:57'
$ws1.Range("D824").Value = '27_EL RETIRO'
$ws1.Range("E824").Value = 55
$ws1.Range("F824").Value = 'LP1912'
$ws1.Range("G824").Value = '31/12/2025'
$ws1.Range("B825").Value = '10:02:06'
$ws1.Range("C825").Value = '11:01'
$ws1.Range("D825").Value = '17_ROMERO'
$ws1.Range("E825").Value = 59
$ws1.Range("F825").Value = 'LP1912'
$ws1.Range("G825").Value = '31/12/2025'
$ws1.Range("B826").Value = '10:02:06'
$ws1.Range("C826").Value = '11:03'
$ws1.Range("D826").Value = '23_HERNANDEZ'
$ws1.Range("E826").Value = 61
$ws1.Range("F826").Value = 'LP1912'
$ws1.Range("G826").Value = '31/12/2025'
$ws1.Range("B827").Value = '10:02:06'
$ws1.Range("C827").Value = '11:05'
$ws1.Range("D827").Value = '14_ABASTO'
$ws1.Range("E827").Value = 63
$ws1.Range("F827").Value = 'LP1912'
$ws1.Range("G827").Value = '31/12/2025'
$ws1.Range("B828").Value = '10:02:06'
$ws1.Range("C828").Value = '11:09'
$ws1.Range("D828").Value = '16_SANTA ANA'
$ws1.Range("E828").Value = 67
$ws1.Range("F828").Value = 'LP1912'
$ws1.Range("G828").Value = '31/12/2025'
$ws1.Range("B829").Value = '10:02:06'
$ws1.Range("C829").Value = '11:11'
$ws1.Range("D829").Value = '15_ABASTO'
$ws1.Range("E829").Value = 69
$ws1.Range("F829").Value = 'LP1912'
$ws1.Range("G829").Value = '31/12/2025'
$ws1.Range("B830").Value = '10:02:06'
$ws1.Range("C830").Value = '11:23'
$ws1.Range("D830").Value = '14_ABASTO'
$ws1.Range("E830").Value = 81
$ws1.Range("F830").Value = 'LP1912'
$ws1.Range("G830").Value = '31/12/2025'
$ws1.Range("B831").Value = '10:02:06'
$ws1.Range("C831").Value = '11:30'
$ws1.Range("D831").Value = '16_P MOR-SANTA ANA'
$ws1.Range("E831").Value = 88
$ws1.Range("F831").Value = 'LP1912'
$ws1.Range("G831").Value = '31/12/2025'
$ws1.Range("B832").Value = '10:02:06'
$ws1.Range("C832").Value = '11:40'
$ws1.Range("D832").Value = '215A_EL PATO'
$ws1.Range("E832").Value = 98
$ws1.Range("F832").Value = 'LP1912'
$ws1.Range("G832").Value = '31/12/2025'

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 31/12/2025 10:02:17'
$ws2.Range("B62").Value = '31/12/2025'
$ws2.Range("C62").Value = '10:02:06'
$ws2.Range("D62").Value = '11:40'
$ws2.Range("E62").Value = '215A_EL PATO'
$ws2.Range("F62").Value = 98
$ws2.Range("G62").Value = 'LP1912'

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 31/12/2025 10:02:17'
$ws3.Range("A3").Value = 'Total filas: 100'
$ws3.Range("B98").Value = '31/12/2025'
$ws3.Range("C98").Value = '10:02:12'
$ws3.Range("D98").Value = '10:09'
$ws3.Range("E98").Value = '215C_LA PLATA'
$ws3.Range("F98").Value = 7
$ws3.Range("G98").Value = 'L6203'
$ws3.Range("B99").Value = '31/12/2025'
$ws3.Range("C99").Value = '10:02:17'
$ws3.Range("D99").Value = '10:23'
$ws3.Range("E99").Value = '215A_LA PLATA'
$ws3.Range("F99").Value = 21
$ws3.Range("G99").Value = 'L6173'
$ws3.Range("B100").Value = '31/12/2025'
$ws3.Range("C100").Value = '10:02:17'
$ws3.Range("D100").Value = '10:31'
$ws3.Range("E100").Value = '215B_LP-P MOR-1 Y 57'
$ws3.Range("F100").Value = 29
$ws3.Range("G100").Value = 'L6173'
$ws3.Range("B101").Value = '31/12/2025'
$ws3.Range("C101").Value = '10:02:17'
$ws3.Range("D101").Value = '11:10'
$ws3.Range("E101").Value = '215A_LA PLATA'
$ws3.Range("F101").Value = 68
$ws3.Range("G101").Value = 'L6173'
